# Update "想去人数" (F column) figures across sheets, per the output
# regenerated at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 150
$ws1.Range("F4").Value  = 1782
$ws1.Range("F5").Value  = 3329
$ws1.Range("F6").Value  = 1068
$ws1.Range("F7").Value  = 2212
$ws1.Range("F8").Value  = 2127
$ws1.Range("F9").Value  = 1112
$ws1.Range("F10").Value = 608
$ws1.Range("F11").Value = 22
$ws1.Range("F12").Value = 1675
$ws1.Range("F13").Value = 398
$ws1.Range("F17").Value = 214
$ws1.Range("F19").Value = 635
$ws1.Range("F22").Value = 12265
$ws1.Range("F23").Value = 12318
$ws1.Range("F24").Value = 911
$ws1.Range("F27").Value = 39
$ws1.Range("F28").Value = 24
$ws1.Range("F29").Value = 372
$ws1.Range("F33").Value = 201
$ws1.Range("F34").Value = 588

# Sheet "演出" (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 9
$ws2.Range("F7").Value = 35

# Sheet "全部类型" (All Types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value  = 150
$ws4.Range("F5").Value  = 1782
$ws4.Range("F6").Value  = 3329
$ws4.Range("F7").Value  = 1068
$ws4.Range("F8").Value  = 2212
$ws4.Range("F9").Value  = 2127
$ws4.Range("F10").Value = 1112
$ws4.Range("F11").Value = 608
$ws4.Range("F12").Value = 22
$ws4.Range("F13").Value = 1675
$ws4.Range("F14").Value = 398
$ws4.Range("F18").Value = 9
$ws4.Range("F21").Value = 214
$ws4.Range("F23").Value = 635
$ws4.Range("F26").Value = 12265
$ws4.Range("F27").Value = 12318
$ws4.Range("F28").Value = 911
$ws4.Range("F31").Value = 39
$ws4.Range("F32").Value = 24
$ws4.Range("F33").Value = 372
$ws4.Range("F39").Value = 201
$ws4.Range("F40").Value = 588
$ws4.Range("F41").Value = 35

$wb.Save()
